$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (D) and Volume(1h) (E) values
$updates = @(
    @{ Row = 2;  D = "306.61";     E = "-4.74%" },
    @{ Row = 3;  D = "40.16";      E = "-6.67%" },
    @{ Row = 4;  D = "5.100";      E = "-1.99%" },
    @{ Row = 5;  D = "0.07697";    E = "-6.05%" },
    @{ Row = 6;  D = "4.267";      E = "-1.36%" },
    @{ Row = 7;  D = "1.618";      E = "-12.34%" },
    @{ Row = 8;  D = "0.8778";     E = "-6.58%" },
    @{ Row = 9;  D = "0.09935";    E = "-10.69%" },
    @{ Row = 10; D = "0.1734";     E = "-6.67%" },
    @{ Row = 11; D = "0.08909";    E = "-4.78%" },
    @{ Row = 12; D = "0.04404";    E = "-4.24%" },
    @{ Row = 13; D = "0.1056";     E = "-0.25%" },
    @{ Row = 14; D = "0.001257";   E = "-3.56%" },
    @{ Row = 15; D = "0.005928";   E = "3.26%" },
    @{ Row = 16; D = "3.356";      E = "-0.05%" },
    @{ Row = 17; D = "2.428";      E = "-3.34%" },
    @{ Row = 18; D = $null;        E = "-1.96%" },
    @{ Row = 19; D = "6.994";      E = "-5.74%" },
    @{ Row = 21; D = "0.3136";     E = "19.64%" },
    @{ Row = 22; D = "0.04144";    E = "0.31%" },
    @{ Row = 23; D = "0.001199";   E = "-3.82%" },
    @{ Row = 24; D = "0.004064";   E = "-5.56%" },
    @{ Row = 25; D = $null;        E = "11.24%" },
    @{ Row = 26; D = $null;        E = "0.04%" },
    @{ Row = 38; D = "0.02344";    E = "-13.67%" },
    @{ Row = 39; D = "0.05145";    E = "-6.77%" },
    @{ Row = 40; D = "0.007963";   E = "0.16%" },
    @{ Row = 41; D = "0.1321";     E = "-5.24%" },
    @{ Row = 42; D = "0.006365";   E = "-2.81%" },
    @{ Row = 43; D = "0.001942";   E = "-8.14%" },
    @{ Row = 44; D = "0.008594";   E = "3.78%" },
    @{ Row = 45; D = "0.3053";     E = "-4.93%" },
    @{ Row = 46; D = "0.00006518"; E = "-6.65%" },
    @{ Row = 47; D = $null;        E = "0.27%" },
    @{ Row = 48; D = "0.007000";   E = "98.22%" },
    @{ Row = 49; D = "0.003612";   E = "4.50%" },
    @{ Row = 50; D = $null;        E = "0.27%" },
    @{ Row = 51; D = $null;        E = "0.27%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }
    $cellE = $ws.Cells.Item($u.Row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
    $cellE.Style = "Normal"
}
